# Refresh the "cryptos" price/volume snapshot (GitHub Actions update).
# Price (D) and Volume 1h % (E) columns are refreshed for every coin row;
# rows 37-40 also get re-ranked (Hedera/dogwifhat and Mantle/Kaspa swap
# places). Some Price values look numeric to Excel (e.g. "596.61"), so
# those are entered with a leading apostrophe to force literal text,
# matching the original inline-string cell content, then the cell style
# is reset to "Normal" so no stray quote-prefix style sticks around.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.718.20"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "3.773.93"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("E4").Value = "  -0.04%  "

$c = $ws.Range("D5")
$c.Value = "'596.61"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "

$c = $ws.Range("D6")
$c.Value = "'168.72"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").Value = "3.771.75"
$ws.Range("E7").Value = "  -1.87%  "

$ws.Range("E8").Value = "  +0.01%  "

$c = $ws.Range("D9")
$c.Value = "'0.525"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("E10").Value = "  -0.52%  "

$c = $ws.Range("D11")
$c.Value = "'6.54"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "

$c = $ws.Range("D12")
$c.Value = "'0.453"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.12%  "

$c = $ws.Range("D13")
$c.Value = "'0.0000280"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.63%  "

$c = $ws.Range("D14")
$c.Value = "'36.32"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("D15").Value = "4.410.31"
$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("D16").Value = "3.780.53"
$ws.Range("E16").Value = "  -1.85%  "

$c = $ws.Range("D17")
$c.Value = "'18.61"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "67.712.58"
$ws.Range("E18").Value = "  -1.14%  "

$c = $ws.Range("D19")
$c.Value = "'7.20"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("E20").Value = "  +0.98%  "

$c = $ws.Range("D21")
$c.Value = "'10.54"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.22%  "

$c = $ws.Range("D22")
$c.Value = "'466.78"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "

$c = $ws.Range("D23")
$c.Value = "'0.717"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.98%  "

$c = $ws.Range("D24")
$c.Value = "'0.0000150"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -7.50%  "

$c = $ws.Range("D25")
$c.Value = "'83.67"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

$c = $ws.Range("D26")
$c.Value = "'2.21"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "

$c = $ws.Range("D27")
$c.Value = "'12.10"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "

$c = $ws.Range("D28")
$c.Value = "'10.45"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.33%  "

$ws.Range("E29").Value = "  +0.00%  "

$c = $ws.Range("D30")
$c.Value = "'2.92"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.33%  "

$ws.Range("D31").Value = "3.925.52"
$ws.Range("E31").Value = "  -1.79%  "

$c = $ws.Range("D32")
$c.Value = "'7.62"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.70%  "

$c = $ws.Range("D33")
$c.Value = "'30.48"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.73%  "

$c = $ws.Range("D34")
$c.Value = "'2.22"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.24%  "

$c = $ws.Range("D35")
$c.Value = "'9.15"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.77%  "

$ws.Range("D36").Value = "3.737.63"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D37")
$c.Value = "'3.78"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.73%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D38")
$c.Value = "'0.104"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D39")
$c.Value = "'0.139"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D40")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "

$c = $ws.Range("D41")
$c.Value = "'5.82"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.00%  "

$c = $ws.Range("D42")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "

$c = $ws.Range("D43")
$c.Value = "'0.312"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "

$ws.Range("E44").Value = "  -0.02%  "

$c = $ws.Range("D45")
$c.Value = "'8.65"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "

$c = $ws.Range("D46")
$c.Value = "'1.94"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.02%  "

$c = $ws.Range("D47")
$c.Value = "'45.72"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.78%  "

$c = $ws.Range("D48")
$c.Value = "'398.86"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -5.06%  "

$c = $ws.Range("D49")
$c.Value = "'0.000272"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -8.17%  "

$c = $ws.Range("D50")
$c.Value = "'140.31"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.87%  "

$c = $ws.Range("D51")
$c.Value = "'39.46"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.21%  "
